# Adds violent-crime data for 2024-01-14 across the Citywide Totals,
# By Neighborhood rollup, and affected individual neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 238
$ws.Range("K3").Value = 222
$ws.Range("I4").Value = 1780
$ws.Range("J4").Value = 1763
$ws.Range("K4").Value = 41
$ws.Range("K5").Value = 11
$ws.Range("J6").Value = 11049
$ws.Range("K6").Value = 306
$ws.Range("I7").Value = 26235
$ws.Range("J7").Value = 29206
$ws.Range("K7").Value = 818

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 3
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 13
$ws.Range("I4").Value = 93
$ws.Range("J6").Value = 690
$ws.Range("K6").Value = 17
$ws.Range("I7").Value = 1542
$ws.Range("J7").Value = 1852
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 9
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 4
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J6").Value = 226
$ws.Range("K7").Value = 25
$ws.Range("I8").Value = 1542
$ws.Range("J8").Value = 1852
$ws.Range("K8").Value = 55
$ws.Range("K9").Value = 4
$ws.Range("K10").Value = 2
$ws.Range("K11").Value = 32
$ws.Range("K15").Value = 4
$ws.Range("K19").Value = 17
$ws.Range("K20").Value = 27
$ws.Range("J27").Value = 179
$ws.Range("K27").Value = 3
$ws.Range("K29").Value = 36
$ws.Range("K31").Value = 9
$ws.Range("K33").Value = 40
$ws.Range("K37").Value = 19
$ws.Range("K41").Value = 11
$ws.Range("K42").Value = 26
$ws.Range("J43").Value = 243
$ws.Range("K51").Value = 16
$ws.Range("K52").Value = 23
$ws.Range("K53").Value = 9
$ws.Range("K54").Value = 12
$ws.Range("K55").Value = 10
$ws.Range("K57").Value = 3
$ws.Range("K65").Value = 18
$ws.Range("K67").Value = 32
$ws.Range("K76").Value = 11
$ws.Range("K83").Value = 12
$ws.Range("K84").Value = 6
$ws.Range("K85").Value = 40
$ws.Range("K89").Value = 8
$ws.Range("K93").Value = 4
$ws.Range("K95").Value = 21
$ws.Range("K96").Value = 9
$ws.Range("K97").Value = 8
$ws.Range("K99").Value = 19
$ws.Range("I101").Value = 26235
$ws.Range("J101").Value = 29206
$ws.Range("K101").Value = 818

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 12
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K2").Value = 2
$ws.Range("K7").Value = 6

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 13
$ws.Range("K3").Value = 11
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K4").Value = 1
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 11

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 226

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 11

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 7
$ws.Range("K3").Value = 9
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 1
$ws.Range("K7").Value = 2

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K3").Value = 5
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 5
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 4

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 4

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 7
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K2").Value = 2
$ws.Range("K7").Value = 4

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 8

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 8

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J2").Value = 48
$ws.Range("K2").Value = 1
$ws.Range("J7").Value = 179
$ws.Range("K7").Value = 3

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 3

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 145
$ws.Range("J7").Value = 243

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 16
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 5
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 23

